# Regenerate orders with updated distance/sizes.
#
# The workbook encodes experiment trial metadata as text like
# "Face04_D64_S25" / "Face04_D64_S25_l.png" / "D64" / "S30" spread across
# several columns (Condition, Filename_Left, Filename_Right, Distance,
# Size). This relabels the distance codes D64/D80/D51 -> D69/D86/D55 and
# the size code S30 -> S31 everywhere those tokens occur, leaving S25/S20
# and everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow = $firstRow + $used.Rows.Count - 1
$lastCol = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -is [string]) {
            $nv = $v.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value = $nv
            }
        }
    }
}
